$d = $word.ActiveDocument

# Locate the start of the paragraph text ("Vous allez participer")
$startRng = $d.Content
$startRng.Find.Text = "Vous allez participer"
$startRng.Find.Forward = $true
$startRng.Find.Wrap = 0
$startRng.Find.Execute() | Out-Null
$start = $startRng.Start

# Locate the end of the paragraph text ("ciel nocturne.")
$endRng = $d.Content
$endRng.Find.Text = "qualité du ciel nocturne."
$endRng.Find.Forward = $true
$endRng.Find.Wrap = 0
$endRng.Find.Execute() | Out-Null
$end = $endRng.End

# Build the replacement range spanning the full run sequence
$full = $d.Range($start, $end)

# Replace the many small runs with a single run containing the full,
# updated paragraph text (naming the constellation "du Taureau" instead
# of "Persee"), with no explicit run formatting.
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Vous allez participer à une campagne mondiale d’observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation du Taureau dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l’éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xmlFrag)
